$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 99, shifting existing rows
# 99-101 down to 100-102 (mirrors the diff: a new weekly record is
# inserted before the former row 99, pushing the rest down by one).
$ws.Rows(99).Insert()

# Populate the newly inserted row 99 with the new record's data.
$ws.Range("A99").Value = 9
$ws.Range("B99").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C99").Value = "Metropolitana"
$ws.Range("D99").Value = 44595
$ws.Range("E99").Value = 13
$ws.Range("F99").Value = 100112022
$ws.Range("G99").Value = "Arveja Verde"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 35
$ws.Range("K99").Value = 24000
$ws.Range("L99").Value = 25000
$ws.Range("M99").Value = 24571
$ws.Range("N99").Value = '$/saco 25 kilos'
$ws.Range("O99").Value = "Carahue"
$ws.Range("P99").Value = 983
$ws.Range("Q99").Value = 25
$ws.Range("R99").Value = "Hortaliza"
